# Commit: "revisions for signs + new overview"
# - Flip the sign of the projected cash-flow values in the Cashflow sheet
#   (rows 2-5, columns F:BN) so they are stored as negative outflows.
# - Remove the now-obsolete helper row 17 (leftover formatting-only row)
#   from the Cashflow sheet.
# - Leave the selection on the Cashflow sheet parked at F15 (new "overview"
#   cell) to match where the author ended up after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cashflow")

# --- 1. Negate the cash-flow projection values (rows 2-5, cols F..BN) ---
for ($r = 2; $r -le 5; $r++) {
    for ($c = 6; $c -le 66; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cur = $cell.Value2
        if ($cur -ne $null) {
            $cell.Value = -1 * $cur
        }
    }
}

# --- 2. Delete the stray row 17 (style-only placeholder row) ---
$ws.Rows.Item(17).Delete()

# --- 3. Park the selection at F15 on the Cashflow sheet ---
[void]$ws.Activate()
[void]$ws.Range("F15").Select()
